$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 23:16"

# Update country rows with refreshed COVID-19 data (values + reordering for Barein / Costa Rica)
# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 3407790
$ws.Cells.Item(4, 3).Value = 52144
$ws.Cells.Item(4, 4).Value = 1512537
$ws.Cells.Item(4, 5).Value = 1757520
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 331
$ws.Cells.Item(4, 8).Value = 137733

# Row 5: Brasil
$ws.Cells.Item(5, 1).Value = "Brasil"
$ws.Cells.Item(5, 2).Value = 1864681
$ws.Cells.Item(5, 3).Value = 23869
$ws.Cells.Item(5, 4).Value = 1213512
$ws.Cells.Item(5, 5).Value = 579069
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 608
$ws.Cells.Item(5, 8).Value = 72100

# Row 27: Egipto
$ws.Cells.Item(27, 1).Value = "Egipto"
$ws.Cells.Item(27, 2).Value = 82070
$ws.Cells.Item(27, 3).Value = 912
$ws.Cells.Item(27, 4).Value = 24419
$ws.Cells.Item(27, 5).Value = 53793
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 89
$ws.Cells.Item(27, 8).Value = 3858

# Row 49: Barein
$ws.Cells.Item(49, 1).Value = "Barein"
$ws.Cells.Item(49, 2).Value = 32941
$ws.Cells.Item(49, 3).Value = 471
$ws.Cells.Item(49, 4).Value = 28425
$ws.Cells.Item(49, 5).Value = 4408
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 4
$ws.Cells.Item(49, 8).Value = 108

# Row 50: Suiza
$ws.Cells.Item(50, 1).Value = "Suiza"
$ws.Cells.Item(50, 2).Value = 32883
$ws.Cells.Item(50, 3).Value = 66
$ws.Cells.Item(50, 4).Value = 29500
$ws.Cells.Item(50, 5).Value = 1415
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 1968

# Row 51: Rumania
$ws.Cells.Item(51, 1).Value = "Rumania"
$ws.Cells.Item(51, 2).Value = 32535
$ws.Cells.Item(51, 3).Value = 456
$ws.Cells.Item(51, 4).Value = 21545
$ws.Cells.Item(51, 5).Value = 9106
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 13
$ws.Cells.Item(51, 8).Value = 1884

# Row 71: Costa de Marfil
$ws.Cells.Item(71, 1).Value = "Costa de Marfil"
$ws.Cells.Item(71, 2).Value = 12766
$ws.Cells.Item(71, 3).Value = 323
$ws.Cells.Item(71, 4).Value = 6654
$ws.Cells.Item(71, 5).Value = 6028
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 2
$ws.Cells.Item(71, 8).Value = 84

# Row 83: Costa Rica
$ws.Cells.Item(83, 1).Value = "Costa Rica"
$ws.Cells.Item(83, 2).Value = 7596
$ws.Cells.Item(83, 3).Value = 365
$ws.Cells.Item(83, 4).Value = 2239
$ws.Cells.Item(83, 5).Value = 5327
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 2
$ws.Cells.Item(83, 8).Value = 30

# Row 84: Etiopia
$ws.Cells.Item(84, 1).Value = "Etiopia"
$ws.Cells.Item(84, 2).Value = 7560
$ws.Cells.Item(84, 3).Value = 158
$ws.Cells.Item(84, 4).Value = 2430
$ws.Cells.Item(84, 5).Value = 5003
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 3
$ws.Cells.Item(84, 8).Value = 127

# Row 85: Finlandia
$ws.Cells.Item(85, 1).Value = "Finlandia"
$ws.Cells.Item(85, 2).Value = 7294
$ws.Cells.Item(85, 3).Value = 3
$ws.Cells.Item(85, 4).Value = 6800
$ws.Cells.Item(85, 5).Value = 165
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 329

# Row 91: Guinea
$ws.Cells.Item(91, 1).Value = "Guinea"
$ws.Cells.Item(91, 2).Value = 6141
$ws.Cells.Item(91, 3).Value = 97
$ws.Cells.Item(91, 4).Value = 4862
$ws.Cells.Item(91, 5).Value = 1242
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 37

# Row 132: Ruanda
$ws.Cells.Item(132, 1).Value = "Ruanda"
$ws.Cells.Item(132, 2).Value = 1337
$ws.Cells.Item(132, 3).Value = 38
$ws.Cells.Item(132, 4).Value = 684
$ws.Cells.Item(132, 5).Value = 649
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 4

# Row 146: Republica del Chad
$ws.Cells.Item(146, 1).Value = "Republica del Chad"
$ws.Cells.Item(146, 2).Value = 880
$ws.Cells.Item(146, 3).Value = 6
$ws.Cells.Item(146, 4).Value = 790
$ws.Cells.Item(146, 5).Value = 15
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 1
$ws.Cells.Item(146, 8).Value = 75

# Row 182: Monaco
$ws.Cells.Item(182, 1).Value = "Monaco"
$ws.Cells.Item(182, 2).Value = 109
$ws.Cells.Item(182, 3).Value = 0
$ws.Cells.Item(182, 4).Value = 97
$ws.Cells.Item(182, 5).Value = 8
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 4

